$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Khách hàng: Đặt lịch, nhận thông báo." -> split the sentence into its
#    own run, colored red and highlighted yellow (the leading ": " stays
#    plain). The following run (a <w:br/>) also becomes red.
# ---------------------------------------------------------------------------
$f1 = $d.Content
$f1.Find.ClearFormatting()
$f1.Find.Replacement.ClearFormatting()
$f1.Find.Text = "Đặt lịch, nhận thông báo."
$f1.Find.Replacement.Text = "Đặt lịch, nhận thông báo."
$f1.Find.Replacement.Highlight = $true
$f1.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $true, $null, 2) | Out-Null

$f2 = $d.Content
$f2.Find.Execute("Đặt lịch, nhận thông báo.", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$f2.Font.Color = 255

$brRange = $d.Range($f2.End, $f2.End + 1)
$brRange.Font.Color = 255

# ---------------------------------------------------------------------------
# 2) Add a collapsed "_GoBack" bookmark right after
#    "...: Quản lý nhân viên, dịch vụ, tài chính." (end of that paragraph).
#    A directly-collapsed Range placed at the paragraph's trailing boundary
#    is mis-seated by this host, so stage it through a throwaway character:
#    insert a marker, wrap the bookmark around it, then delete the marker.
# ---------------------------------------------------------------------------
$f3 = $d.Content
$f3.Find.Execute("Quản lý nhân viên, dịch vụ, tài chính.", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$pos = $f3.End
$marker = $d.Range($pos, $pos)
$marker.InsertAfter("X")
$markerRange = $d.Range($pos, $pos + 1)
$d.Bookmarks.Add("_GoBack", $markerRange) | Out-Null
$cleanup = $d.Range($pos, $pos + 1)
$cleanup.Text = ""

# ---------------------------------------------------------------------------
# 3) "2. Mô hình dữ liệ" + old "_GoBack" bookmark + "u (Database - ERD)" ->
#    merge into a single run "2. Mô hình dữ liệu (Database - ERD)" and
#    drop the stale bookmark that used to sit mid-word.
# ---------------------------------------------------------------------------
$old = $d.Bookmarks("_GoBack")
$old.Delete()

$f4 = $d.Content
$f4.Find.ClearFormatting()
$f4.Find.Replacement.ClearFormatting()
$f4.Find.Text = "2. Mô hình dữ liệu (Database - ERD)"
$f4.Find.Replacement.Text = "2. Mô hình dữ liệu (Database - ERD)"
$f4.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null
